$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 2, shifting existing rows 2-70 (A:AA only) down
# to 3-70. Restrict the insert to the used range (A:AA) so Excel doesn't
# stamp formatting across the entire 16384-column row.
$ws.Range("A2:AA2").Insert(-4121)  # xlShiftDown

# The new row 2 (A2:AA2) inherits the generic style by default. Columns A/B
# need the date-number-format style used by the rest of column A/B, so copy
# that formatting down from row 3 (which holds what used to be row 2).
$ws.Range("A3:B3").Copy()
$ws.Range("A2:B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Make sure the new row explicitly carries the same row height as the rest
# of the sheet (matches the ht/customHeight attributes on other rows).
$ws.Rows.Item(2).RowHeight = 15.75

# Populate the new row 2 with the new data values.
$ws.Range("A2").Value = 43927
$ws.Range("B2").Value = 43930
$ws.Range("C2").Value = 257.83999999999997
$ws.Range("D2").Value = 281.2
$ws.Range("E2").Value = 248.17
$ws.Range("F2").Value = 278.2
$ws.Range("G2").Value = 259.94099999999997
$ws.Range("H2").Value = 236.43899999999999
